$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the greeting text in E8 from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Set active cell / selection to E8 (as reflected in the sheetView)
$ws.Activate()
$ws.Range("E8").Select()
